$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded ahead of the existing data set.
# Insert a new row at position 49, which pushes the existing rows 49-119
# down to 50-120 (dimension grows from R119 to R120).
$ws.Rows(49).Insert()

# Populate the newly inserted row 49 with the new observation. All
# descriptive columns mirror the row immediately below (same market,
# region, product, etc.) while the date and price columns carry the new
# figures.
$ws.Cells.Item(49, 1).Value = 9
$ws.Cells.Item(49, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(49, 3).Value = "Metropolitana"
$ws.Cells.Item(49, 4).Value = 45210
$ws.Cells.Item(49, 5).Value = 13
$ws.Cells.Item(49, 6).Value = 100112029
$ws.Cells.Item(49, 7).Value = "Orégano"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 16
$ws.Cells.Item(49, 11).Value = 21000
$ws.Cells.Item(49, 12).Value = 21000
$ws.Cells.Item(49, 13).Value = 21000
$ws.Cells.Item(49, 14).Value = "$/docena de atados"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 7000
$ws.Cells.Item(49, 17).Value = 3
$ws.Cells.Item(49, 18).Value = "Hortaliza"
